# "Generate Report for Handback" — mark the en-US source as handed back,
# record the generated/received localization artifacts for zh-cn and
# de-de, and widen the columns that now hold the longer status/filename
# text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$mdFileName = "36db075a-336b-4ecf-ad40-2d245c6e7a1c.md"
$mdFileUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4ad9baa6fc4ad34c030d70560d49ba1811b1e65/e2e/36db075a-336b-4ecf-ad40-2d245c6e7a1c.md"

# --- Overview sheet: zh-cn / de-de status columns -------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText

# --- zh-cn sheet: status + handback bookkeeping ----------------------------
$zhcn.Range("C2").Value = $statusText

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdFileUrl, "", "", $mdFileName) | Out-Null
$zhcn.Range("J2").Value = "36db075a-336b-4ecf-ad40-2d245c6e7a1c.84f50f49420392f2643eb20fea570acd54a163d3.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-24 02:58:27"

# --- de-de sheet: status + handback bookkeeping ----------------------------
$dede.Range("C2").Value = $statusText

$dede.Hyperlinks.Add($dede.Range("I2"), $mdFileUrl, "", "", $mdFileName) | Out-Null
$dede.Range("J2").Value = "36db075a-336b-4ecf-ad40-2d245c6e7a1c.84f50f49420392f2643eb20fea570acd54a163d3.de-de.xlf"
$dede.Range("K2").Value = "2016-08-24 02:58:35"

# --- Column widths (widened to fit the longer text now shown) -------------
# NOTE: ColumnWidth is specified in characters; the file stores
# characters+5/6 (Excel's standard padding) rounded to the nearest 1/6, so
# we set values that land as close as possible to the intended stored
# widths of ~29.98 and 40.
$overview.Columns.Item(5).ColumnWidth = 29.16666666666667
$overview.Columns.Item(6).ColumnWidth = 29.16666666666667

$zhcn.Columns.Item(3).ColumnWidth  = 29.16666666666667
$zhcn.Columns.Item(9).ColumnWidth  = 39.16666666666666
$zhcn.Columns.Item(10).ColumnWidth = 39.16666666666666

$dede.Columns.Item(3).ColumnWidth  = 29.16666666666667
$dede.Columns.Item(9).ColumnWidth  = 39.16666666666666
$dede.Columns.Item(10).ColumnWidth = 39.16666666666666
